$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CommonTestdata")
$ws2 = $wb.Worksheets.Item("Headers")

# ---------------------------------------------------------------
# Sheet "Headers" -> "Fieldinfo"
# ---------------------------------------------------------------
$ws2.Name = "Fieldinfo"

# ---------------------------------------------------------------
# Sheet1 (CommonTestdata): testdata no longer carries a login/pwd
# ---------------------------------------------------------------
$ws1.Range("D2:E2").ClearContents()
$ws1.Range("D2:E2").Select()

# ---------------------------------------------------------------
# Sheet2 (Fieldinfo): extend the table from 6 to 13 columns and
# add the new "field info" row describing the Forgot-password page
# ---------------------------------------------------------------
$ws2.Activate()

# Bring the header format (grey fill + border, style of A1) onto the
# newly added header cells G1:M1.
$ws2.Range("A1").Copy()
$ws2.Range("G1:M1").PasteSpecial(-4122)

# Bring the plain bordered body format (style of A2) onto the newly
# added body cells H2:M8 (columns H..M, rows 2..8).
$ws2.Range("A2").Copy()
$ws2.Range("H2:M8").PasteSpecial(-4122)

# Columns E:G reuse the bordered "date" look already used by D:E in
# the old layout (style of D2), rows 2,4-8. Column D itself reverts to
# the plain bordered look (it now holds "Header2", not a date).
$ws2.Range("D2").Copy()
$ws2.Range("E2:G2").PasteSpecial(-4122)
$ws2.Range("E4:G8").PasteSpecial(-4122)

$ws2.Range("A2").Copy()
$ws2.Range("D2").PasteSpecial(-4122)
$ws2.Range("D4:D8").PasteSpecial(-4122)
$ws2.Range("D3").PasteSpecial(-4122)
$ws2.Range("E3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Header row text
$ws2.Range("A1").Value = "Tc_Name"
$ws2.Range("B1").Value = "Iteration"
$ws2.Range("C1").Value = "Header1"
$ws2.Range("D1").Value = "Header2"
$ws2.Range("E1").Value = "Text1"
$ws2.Range("F1").Value = "Text2"
$ws2.Range("G1").Value = "Email"
$ws2.Range("H1").Value = "Error1"
$ws2.Range("I1").Value = "Error2"
$ws2.Range("J1").Value = "Rtrpwd"
$ws2.Range("K1").Value = "Bck2Login"
$ws2.Range("L1").Value = "Forgotlink"
$ws2.Range("M1").Value = "Authlink"

# Row2 data
$ws2.Range("A2").Value = "TC_01_Verify_Forgotpwd_link_is_available"
$ws2.Range("B2").Value = 1

# Row3 data
$ws2.Range("A3").Value = "TC_02_Verify_all_fields_Frogotpwd_page"
$ws2.Range("B3").Value = 1
$ws2.Range("D3").Value = "Forgot your password?"
$ws2.Range("E3").Value = "Please enter the email address you used to register."

# C3, F3, G3 get the "date format, no quote-prefix" look: copy the
# date-format style from D2 then overwrite with literal text (a typed
# plain Value assignment never sets quotePrefix, unlike a leading ').
$ws2.Range("D2").Copy()
$ws2.Range("C3").PasteSpecial(-4122)
$ws2.Range("F3").PasteSpecial(-4122)
$ws2.Range("G3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Range("C3").Value = "Authentication"
$ws2.Range("F3").Value = "Email address"
$ws2.Range("G3").Value = "abcd"

$ws2.Range("H3").Value = "There is 1 error"
$ws2.Range("I3").Value = "Invalid email address."
$ws2.Range("J3").Value = "Retrieve Password"
$ws2.Range("K3").Value = "Back to Login"
$ws2.Range("L3").Value = "Forgot your password?"
$ws2.Range("M3").Value = "Authentication"

# Column widths (approximate AutoFit to content, matches the author's
# widened table after the new columns were populated).
$ws2.Range("A1:M8").Columns.AutoFit()

$ws2.Range("G1").Select()

Write-Host "Sheet1:" $ws1.Name
Write-Host "Sheet2:" $ws2.Name
